$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (61) with the 23/4/2019 locale counts, one day after
# the previous last row (60, date 22/4/2019 -> serial 43577).
$row = 61

$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item($row, 1).Value = 43578

$values = @(6, 5, 1, 459, 17, 35, 3, 1, 3, 1, 5, 4, 2, 1, 13, 2, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($row, $i + 2).Value = $values[$i]
}

# Match the view state shown in the diff: scroll so row 38 is the top
# visible row, and select the newly added last cell A61.
$ws.Application.ActiveWindow.ScrollRow = 38
$ws.Range("A61").Select()
